# feat: add Usecase_Scenarios and Planification_initiale in pdf
#
# Journal de Travail: row 9's duration cell becomes a text entry, and six more
# rows (9 more cells of new text "Planification initiale" follow-up + the new
# "Usecase_Scenarios"/"Planification du MCD" entries) are filled in below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de Travail")

# Re-use the existing "date" (B6) and "time, centered" (C6) cell formatting so
# the newly-filled rows look like the rest of the table, instead of minting
# brand-new style/numFmt entries.
$ws.Range("B6").Copy()
$ws.Range("B10:B14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C6").Copy()
$ws.Range("C9:C14").PasteSpecial(-4122)    # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 9: the 0h40 duration is replaced by a plain-text duration, and the
#     activity description is updated ---
$ws.Range("C9").Value = "35 min"
$ws.Range("D9").Value = " Début de la Planification Initiale du projet"

# --- Row 10: continuation of "Planification Initiale du projet" ---
$ws.Range("B10").Value = 46056
$ws.Range("C10").Value = 0.3888888888888889
$ws.Range("D10").Value = "Planification Initiale du projet"
$ws.Range("E10").Value = "Continuation de la planification du projet"

# --- Row 11: finalisation of the planning ---
$ws.Range("B11").Value = 46056
$ws.Range("C11").Value = 0.44444444444444442
$ws.Range("D11").Value = "Finalisation de la planification du projet"
$ws.Range("E11").Value = "Fin de la Planification du projet"

# --- Row 12: Use cases ---
$ws.Range("B12").Value = 46056
$ws.Range("C12").Value = "30 min"
$ws.Range("D12").Value = "Création des Use cases"
$ws.Range("E12").Value = "Aide avec Gemini pour les Use cases"

# --- Row 13: Scenarios ---
$ws.Range("B13").Value = 46056
$ws.Range("C13").Value = "50 min "
$ws.Range("D13").Value = "Scenarios"
$ws.Range("E13").Value = "Les scénarios sont créent dans le même fichier que les use cases"

# --- Row 14: MCD ---
$ws.Range("B14").Value = 46056
$ws.Range("C14").Value = 0.65277777777777779
$ws.Range("D14").Value = "Création du MCD"
$ws.Range("E14").Value = "LE MCD est créer avec draw.io"

# Columns C and D grew wider to fit the new, longer text.
$ws.Columns.Item(3).ColumnWidth = 17.67
$ws.Columns.Item(4).ColumnWidth = 37.67

# The cursor ends up on the next empty row.
$ws.Range("B15").Select()
